$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive updated values for each data row (rows 2-25)
$cols = @("B","C","D","E","I","J","K","L","N")

# Each string: row number followed by the 9 new values (for columns above, in order)
$data = @(
    "2,1.02,1.017235370415475,1.022499601845181,1.018665282083338,1.026651686960011,1.022450639559592,1.025333801797923,1.021510822975832,1.023902636744495",
    "3,1.02,1.018117027998092,1.0231225573549,1.019410289531132,1.026769453774634,1.022968452246334,1.025763938586235,1.022061857090794,1.024421184784646",
    "4,1.02,1.018687530390998,1.023525114953858,1.019892789334367,1.026843600842647,1.023302916375415,1.026041042186972,1.022418181792492,1.024756123891154",
    "5,1.02,1.018927370907646,1.023694219983082,1.020095733600907,1.026874279136895,1.02344338132738,1.026157242189603,1.022567924053104,1.024896788319461",
    "6,1.02,1.018967641224162,1.023722605774838,1.020129814731765,1.026879401212484,1.023466957543746,1.026176735379431,1.022593063078271,1.024920398016759",
    "7,1.02,1.018690735148235,1.023527375057366,1.019895500689426,1.026844012706388,1.023304793841475,1.02604259601548,1.022420182880116,1.024758004023431",
    "8,1.02,1.017533327658648,1.022710242536151,1.018916970789794,1.026691911940271,1.022625759293657,1.025479421066646,1.02169709463668,1.024078005168669",
    "9,1.02,1.01549396967268,1.021266322367195,1.017196042410193,1.026408185920489,1.021424709145123,1.024477724478011,1.020421206504582,1.02287524939116",
    "10,1.02,1.014134569319631,1.020301105914231,1.016051116015023,1.026208530526977,1.020621053323673,1.02380374874295,1.019569537242814,1.022070452286264",
    "11,1.02,1.013545990226507,1.01988256104529,1.015555929771821,1.026119597194223,1.020272375471675,1.02351046138395,1.01920051518807,1.021721279271721",
    "12,1.02,1.013327374485891,1.019727006468786,1.015372083443577,1.026086191479045,1.020142758565951,1.023401304661507,1.019063408591491,1.021591478295125",
    "13,1.02,1.013374267889685,1.019760377425018,1.01541151510243,1.026093373934574,1.020170566448464,1.023424728935804,1.019092820002985,1.02161932566802",
    "14,1.02,1.013527919195985,1.019869704648348,1.015540731174885,1.026116843449401,1.020261663380519,1.023501442869027,1.019189182628572,1.021710551968168",
    "15,1.02,1.013622590030851,1.019937053099774,1.015620357155494,1.026131254535265,1.020317777669422,1.0235486801745,1.019248550171735,1.021766745945798",
    "16,1.02,1.014173632535974,1.020328870870377,1.016083992133628,1.026214380536402,1.020644179496197,1.023823182792292,1.019594022993102,1.022093611300605",
    "17,1.02,1.014519301234118,1.020574488456754,1.016374973023818,1.026265859595559,1.020848738950854,1.023994983381915,1.019810664606907,1.022298461253151",
    "18,1.02,1.014720928704195,1.020717695054238,1.016544752573334,1.026295646999817,1.020967988407267,1.024095051630443,1.019937004404745,1.022417880057471",
    "19,1.02,1.014789679233066,1.020766514914993,1.016602652346913,1.026305763092007,1.021008638055712,1.024129148508595,1.019980078972053,1.02245858743308",
    "20,1.02,1.014482213752447,1.020548141999344,1.016343747794682,1.026260361141917,1.020826798538825,1.023976565289139,1.019787423443144,1.022276489683221",
    "21,1.02,1.013482672465072,1.019837512931754,1.015502677823254,1.026109942522063,1.02023484042264,1.023478858500233,1.019160807203785,1.021683690918612",
    "22,1.02,1.012854272346968,1.01939020176197,1.01497437303826,1.026013216996972,1.01986206082751,1.023164677150981,1.018766624640126,1.021310381933699",
    "23,1.02,1.013187393892253,1.019627377687804,1.015254388598159,1.02606469666741,1.020059734071737,1.023331348992459,1.018975607209363,1.021508335896606",
    "24,1.02,1.014498971974349,1.020560047008715,1.01635785696089,1.026262846396402,1.02083671267253,1.023984888065559,1.019797925197649,1.022286417896133",
    "25,1.02,1.016021167902524,1.021640076804447,1.017640534809344,1.026483391380355,1.021735735667254,1.024737781916939,1.02075124951989,1.023186717606628"
)

foreach ($line in $data) {
    $parts = $line.Split(",")
    $row = [int]$parts[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $val = [double]$parts[$i + 1]
        $cellRef = $cols[$i] + $row
        $ws.Range($cellRef).Value = $val
    }
}
